# AdvancedRepeatForm_import_test_01.xlsx edit:
#  - rename Sheet1 ("Advanced Repeat Form") headers K1:M1 to add a
#    "group_recent_haircuts::" prefix (they now live inside that group)
#  - add a new worksheet "IDSheet" (after group_pets) holding the form's
#    KPI/KC identifiers, stored as text

$wb = $excel.ActiveWorkbook

# --- rename the repeat-group headers on the main sheet ---------------------
$ws1 = $wb.Worksheets.Item("Advanced Repeat Form")
$ws1.Range("K1").Value = "group_recent_haircuts::Last_Haircut"
$ws1.Range("L1").Value = "group_recent_haircuts::Haircut_before_last"
$ws1.Range("M1").Value = "group_recent_haircuts::Three_haircuts_ago"

# --- add the new IDSheet as the last tab ------------------------------------
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$idSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$idSheet.Name = "IDSheet"

$idSheet.Range("A1").Value = "KPI ID"
$idSheet.Range("B1").Value = "afRmyXNV6NeAtcXYs3fNaD"
$idSheet.Range("A2").Value = "KC ID"
$idSheet.Range("B2").Value = "bb97d983300345d88748e629139f3062"

$ws1.Select() | Out-Null
$ws1.Range("K12").Select() | Out-Null
